$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3605800351898267
$ws.Range("C2").Value = 0.03088557933001823
$ws.Range("D2").Value = 0.1582768291400214
$ws.Range("E2").Value = 0.1566086668041855
$ws.Range("F2").Value = 1.732319223695441
$ws.Range("J2").Value = 0.1935279873022253
$ws.Range("K2").Value = 0.3142421109736802
$ws.Range("M2").Value = 0.2122778647614609
$ws.Range("O2").Value = 4.396756074374395
$ws.Range("B3").Value = 0.3283216481628415
$ws.Range("C3").Value = 0.02738043923127975
$ws.Range("D3").Value = 0.1550313151750942
$ws.Range("E3").Value = 0.1554378733550017
$ws.Range("F3").Value = 1.737500715248807
$ws.Range("J3").Value = 0.1934081094815738
$ws.Range("K3").Value = 0.2810876666783031
$ws.Range("M3").Value = 0.2020710694219403
$ws.Range("O3").Value = 4.419837966818903
$ws.Range("B4").Value = 0.3085938214299802
$ws.Range("C4").Value = 0.02521660268801895
$ws.Range("D4").Value = 0.1530991578749621
$ws.Range("E4").Value = 0.1547863530560392
$ws.Range("F4").Value = 1.741507268733137
$ws.Range("J4").Value = 0.1934220207013269
$ws.Range("K4").Value = 0.26075683036089
$ws.Range("M4").Value = 0.1958871293921334
$ws.Range("O4").Value = 4.43613697541997
$ws.Range("B5").Value = 0.3005748792314478
$ws.Range("C5").Value = 0.02433193968309411
$ws.Range("D5").Value = 0.1523271114385096
$ws.Range("E5").Value = 0.1545378302330747
$ws.Range("F5").Value = 1.743347577505148
$ws.Range("J5").Value = 0.193449718883663
$ws.Range("K5").Value = 0.2524788782477003
$ws.Range("M5").Value = 0.1933881755924247
$ws.Range("O5").Value = 4.443313832954075
$ws.Range("B6").Value = 0.2992445799039274
$ws.Range("C6").Value = 0.02418486951992804
$ws.Range("D6").Value = 0.1521998417017301
$ws.Range("E6").Value = 0.1544975897796519
$ws.Range("F6").Value = 1.743665702032885
$ws.Range("J6").Value = 0.193455649320633
$ws.Range("K6").Value = 0.2511047686838026
$ws.Range("M6").Value = 0.1929745024975418
$ws.Range("O6").Value = 4.444537853625192
$ws.Range("B7").Value = 0.3084855923441694
$ws.Range("C7").Value = 0.02520468341344895
$ws.Range("D7").Value = 0.153088683629079
$ws.Range("E7").Value = 0.1547829325966248
$ws.Range("F7").Value = 1.741531247024781
$ws.Range("J7").Value = 0.1934223050231623
$ws.Range("K7").Value = 0.2606451619339509
$ws.Range("M7").Value = 0.1958533421725548
$ws.Range("O7").Value = 4.436231599172203
$ws.Range("B8").Value = 0.3494412181356097
$ws.Range("C8").Value = 0.02967945817371742
$ws.Range("D8").Value = 0.1571452475175477
$ws.Range("E8").Value = 0.1561910222056397
$ws.Range("F8").Value = 1.733934627502101
$ws.Range("J8").Value = 0.1934685023262261
$ws.Range("K8").Value = 0.3028053005196512
$ws.Range("M8").Value = 0.2087414135881644
$ws.Range("O8").Value = 4.404273395979715
$ws.Range("B9").Value = 0.4303657638029961
$ws.Range("C9").Value = 0.03836005946548937
$ws.Range("D9").Value = 0.1655778204648755
$ws.Range("E9").Value = 0.1594851904313579
$ws.Range("F9").Value = 1.725580150082749
$ws.Range("J9").Value = 0.1942527371953773
$ws.Range("K9").Value = 0.385672512603719
$ws.Range("M9").Value = 0.2346683482246448
$ws.Range("O9").Value = 4.358474400888497
$ws.Range("B10").Value = 0.4901782723370616
$ws.Range("C10").Value = 0.04467825869863873
$ws.Range("D10").Value = 0.1720609437847287
$ws.Range("E10").Value = 0.1622287232871464
$ws.Range("F10").Value = 1.723426786165774
$ws.Range("J10").Value = 0.1952510818985118
$ws.Range("K10").Value = 0.446656336400423
$ws.Range("M10").Value = 0.2541097358512729
$ws.Range("O10").Value = 4.335110384855483
$ws.Range("B11").Value = 0.5174633051149158
$ws.Range("C11").Value = 0.04753930817601315
$ws.Range("D11").Value = 0.1750720700495947
$ws.Range("E11").Value = 0.1635467263091535
$ws.Range("F11").Value = 1.723311738136658
$ws.Range("J11").Value = 0.1957967994372325
$ws.Range("K11").Value = 0.4744185936270071
$ws.Range("M11").Value = 0.263038403190734
$ws.Range("O11").Value = 4.326714858245566
$ws.Range("B12").Value = 0.5278059967522495
$ws.Range("C12").Value = 0.04862078184997642
$ws.Range("D12").Value = 0.1762211348655995
$ws.Range("E12").Value = 0.1640558447367759
$ws.Range("F12").Value = 1.723392402451481
$ws.Range("J12").Value = 0.1960165987543476
$ws.Range("K12").Value = 0.4849339962834449
$ws.Range("M12").Value = 0.2664315005847584
$ws.Range("O12").Value = 4.323856747531408
$ws.Range("B13").Value = 0.5255780564370127
$ws.Range("C13").Value = 0.04838795446893585
$ws.Range("D13").Value = 0.1759732724673739
$ws.Range("E13").Value = 0.1639457520059402
$ws.Range("F13").Value = 1.723369506162157
$ws.Range("J13").Value = 0.1959686765853448
$ws.Range("K13").Value = 0.4826692156169372
$ws.Range("M13").Value = 0.2657002051737507
$ws.Range("O13").Value = 4.324458012003362
$ws.Range("B14").Value = 0.5183139981113811
$ws.Range("C14").Value = 0.04762832093601332
$ws.Range("D14").Value = 0.1751664281053706
$ws.Range("E14").Value = 0.1635884112181394
$ws.Range("F14").Value = 1.723315885177428
$ws.Range("J14").Value = 0.1958146190646701
$ws.Range("K14").Value = 0.4752836561013964
$ws.Range("M14").Value = 0.2633173157468462
$ws.Range("O14").Value = 4.326473284535524
$ws.Range("B15").Value = 0.5138658979627166
$ws.Range("C15").Value = 0.04716276870186675
$ws.Range("D15").Value = 0.17467335823018
$ws.Range("E15").Value = 0.1633708332775363
$ws.Range("F15").Value = 1.723299216647789
$ws.Range("J15").Value = 0.1957219660661167
$ws.Range("K15").Value = 0.470760090485328
$ws.Range("M15").Value = 0.2618592854425614
$ws.Range("O15").Value = 4.327749513144909
$ws.Range("B16").Value = 0.4883966191225966
$ws.Range("C16").Value = 0.04449101340861716
$ws.Range("D16").Value = 0.1718653981264708
$ws.Range("E16").Value = 0.1621439929404858
$ws.Range("F16").Value = 1.723451693864021
$ws.Range("J16").Value = 0.1952172588763403
$ws.Range("K16").Value = 0.4448423758111915
$ws.Range("M16").Value = 0.2535279157281494
$ws.Range("O16").Value = 4.335703984193259
$ws.Range("B17").Value = 0.4727911635946214
$ws.Range("C17").Value = 0.0428485753289749
$ws.Range("D17").Value = 0.1701586026832587
$ws.Range("E17").Value = 0.1614092548985369
$ws.Range("F17").Value = 1.723766597907826
$ws.Range("J17").Value = 0.1949310759262204
$ws.Range("K17").Value = 0.4289475709117312
$ws.Range("M17").Value = 0.2484384567399971
$ws.Range("O17").Value = 4.341155694520097
$ws.Range("B18").Value = 0.4638225031433478
$ws.Range("C18").Value = 0.04190265572017893
$ws.Range("D18").Value = 0.1691827326036446
$ws.Range("E18").Value = 0.1609932400281231
$ws.Range("F18").Value = 1.724029101111853
$ws.Range("J18").Value = 0.1947750907873598
$ws.Range("K18").Value = 0.4198072501792183
$ws.Range("M18").Value = 0.2455191155556591
$ws.Range("O18").Value = 4.344501543058755
$ws.Range("B19").Value = 0.4607871177724974
$ws.Range("C19").Value = 0.04158217357743865
$ws.Range("D19").Value = 0.1688533243429049
$ws.Range("E19").Value = 0.1608535170959513
$ws.Range("F19").Value = 1.724131959725099
$ws.Range("J19").Value = 0.1947237579262193
$ws.Range("K19").Value = 0.4167128448162885
$ws.Range("M19").Value = 0.2445320527966572
$ws.Range("O19").Value = 4.345670486379447
$ws.Range("B20").Value = 0.4744516509289838
$ws.Range("C20").Value = 0.04302354364648409
$ws.Range("D20").Value = 0.1703396908594499
$ws.Range("E20").Value = 0.161486787578017
$ws.Range("F20").Value = 1.723724654072825
$ws.Range("J20").Value = 0.1949606486189666
$ws.Range("K20").Value = 0.4306394017466459
$ws.Range("M20").Value = 0.2489794135828944
$ws.Range("O20").Value = 4.340553598888562
$ws.Range("B21").Value = 0.5204473486347752
$ws.Range("C21").Value = 0.04785149690881951
$ws.Range("D21").Value = 0.1754031792677893
$ws.Range("E21").Value = 0.1636930992705992
$ws.Range("F21").Value = 1.72332826408045
$ws.Range("J21").Value = 0.195859512813918
$ws.Range("K21").Value = 0.4774529117315467
$ws.Range("M21").Value = 0.2640169036979572
$ws.Range("O21").Value = 4.325872635714461
$ws.Range("B22").Value = 0.5505688328310896
$ws.Range("C22").Value = 0.05099548880559723
$ws.Range("D22").Value = 0.1787638113205361
$ws.Range("E22").Value = 0.1651934328514457
$ws.Range("F22").Value = 1.723793269801362
$ws.Range("J22").Value = 0.1965235998339949
$ws.Range("K22").Value = 0.5080622301992435
$ws.Range("M22").Value = 0.2739146664318355
$ws.Range("O22").Value = 4.318149356316781
$ws.Range("B23").Value = 0.5344870431331117
$ws.Range("C23").Value = 0.04931853833156197
$ws.Range("D23").Value = 0.1769655091253668
$ws.Range("E23").Value = 0.1643873480144364
$ws.Range("F23").Value = 1.723478864147708
$ws.Range("J23").Value = 0.1961621590051337
$ws.Range("K23").Value = 0.4917243449738464
$ws.Range("M23").Value = 0.2686257057968575
$ws.Range("O23").Value = 4.322100164057531
$ws.Range("B24").Value = 0.4737009348388028
$ws.Range("C24").Value = 0.04294444563628019
$ws.Range("D24").Value = 0.1702578041017659
$ws.Range("E24").Value = 0.1614517151286634
$ws.Range("F24").Value = 1.723743363127298
$ws.Range("J24").Value = 0.1949472521857984
$ws.Range("K24").Value = 0.4298745317984753
$ws.Range("M24").Value = 0.2487348262474569
$ws.Range("O24").Value = 4.340825147233147
$ws.Range("B25").Value = 0.4084097131299131
$ws.Range("C25").Value = 0.03602203066499499
$ws.Range("D25").Value = 0.1632458097352583
$ws.Range("E25").Value = 0.1585371411837073
$ws.Range("F25").Value = 1.727140232460336
$ws.Range("J25").Value = 0.1939663935490472
$ws.Range("K25").Value = 0.3632358490208389
$ws.Range("M25").Value = 0.2275850375963984
$ws.Range("O25").Value = 4.369058140513516
